$wb = $excel.ActiveWorkbook

# --- Style source sheets (existing sheets with the formatting we need) ---
$detailStyleSrc = $wb.Worksheets.Item("2021-Q4")   # has the 8-col fund-holding header/index style
$summaryStyleSrc = $wb.Worksheets.Item("2021-Q1")  # has the 4-col quarter-summary header/index style

# --- Step 1: the existing "总计" sheet becomes the new "2022-Q1" detail sheet ---
$detail = $wb.Worksheets.Item("总计")
$detail.Name = "2022-Q1"
$detail.Outline.SummaryRow = 1
$detail.Outline.SummaryColumn = 1
$detail.Cells.Clear()

# Copy header-row formatting (bold / centered / bordered) and column-A index formatting
$detailStyleSrc.Range("B1:H1").Copy()
$detail.Range("B1:H1").PasteSpecial(-4122)
$detailStyleSrc.Range("A2").Copy()
$detail.Range("A2:A4").PasteSpecial(-4122)

# Header row
$detail.Range("B1").Value = "基金代码"
$detail.Range("C1").Value = "基金名称"
$detail.Range("D1").Value = "基金规模"
$detail.Range("E1").Value = "股票总仓位"
$detail.Range("F1").Value = "仓位占比"
$detail.Range("G1").Value = "持有市值(亿元)"
$detail.Range("H1").Value = "仓位排名"

# Index column (0-based row counter)
$detail.Range("A2").Value = 0
$detail.Range("A3").Value = 1
$detail.Range("A4").Value = 2

# Columns B-G hold text (codes / names / formatted numeric strings)
$detail.Range("B2:G4").NumberFormat = "@"

$detail.Range("B2").Value = "005571"
$detail.Range("C2").Value = "中银证券新能源灵活配置混合A"
$detail.Range("D2").Value = "0.91"
$detail.Range("E2").Value = "90.25"
$detail.Range("F2").Value = "8.92"
$detail.Range("G2").Value = "0.0812"
$detail.Range("H2").Value = 2

$detail.Range("B3").Value = "005572"
$detail.Range("C3").Value = "中银证券新能源灵活配置混合C"
$detail.Range("D3").Value = "0.28"
$detail.Range("E3").Value = "90.25"
$detail.Range("F3").Value = "8.92"
$detail.Range("G3").Value = "0.0250"
$detail.Range("H3").Value = 2

$detail.Range("B4").Value = "003981"
$detail.Range("C4").Value = "中银证券瑞益灵活配置混合C"
$detail.Range("D4").Value = "0.21"
$detail.Range("E4").Value = "89.21"
$detail.Range("F4").Value = "6.06"
$detail.Range("G4").Value = "0.0127"
$detail.Range("H4").Value = 1

# --- Step 2: add a brand-new "总计" sheet after "2022-Q1", with the refreshed quarterly roll-up ---
$summary = $wb.Worksheets.Add($null, $detail)
$summary.Name = "总计"
$summary.Outline.SummaryRow = 1
$summary.Outline.SummaryColumn = 1

$summaryStyleSrc.Range("B1:D1").Copy()
$summary.Range("B1:D1").PasteSpecial(-4122)
$summaryStyleSrc.Range("A2").Copy()
$summary.Range("A2:A7").PasteSpecial(-4122)

$summary.Range("B1").Value = "日期"
$summary.Range("C1").Value = "持有数量(只)"
$summary.Range("D1").Value = "持有市值(亿元)"

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q1"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.12

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 10
$summary.Range("D3").Value = 2.43

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 39
$summary.Range("D4").Value = 5.83

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 27
$summary.Range("D5").Value = 13.47

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 19
$summary.Range("D6").Value = 9.140000000000001

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2020-Q4"
$summary.Range("C7").Value = 13
$summary.Range("D7").Value = 5.86

# Restore the originally-active tab
$wb.Worksheets.Item("2020-Q4").Activate()
